$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all cells being updated so numeric-looking strings
# (e.g. "582.81", "1.00") are preserved exactly as text, matching the source data.
$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "E10", "E11", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "E34", "E35", "D36", "E36", "E37", "E38", "D39", "E39", "D40", "E40", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "B47", "C47", "D47", "E47", "B48", "C48", "D48", "E48", "B49", "C49", "D49", "E49", "D50", "E50", "E51")
foreach ($ref in $cells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "60.542.25"
$ws.Range("E2").Value = "  +2.12%  "
$ws.Range("D3").Value = "2.626.91"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "582.81"
$ws.Range("E5").Value = "  +5.77%  "
$ws.Range("D6").Value = "143.70"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("D9").Value = "2.627.73"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  -2.27%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("E12").Value = "  -3.88%  "
$ws.Range("D13").Value = "0.374"
$ws.Range("E13").Value = "  +6.04%  "
$ws.Range("D14").Value = "3.084.93"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "25.08"
$ws.Range("E15").Value = "  +8.60%  "
$ws.Range("D16").Value = "60.552.28"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  +3.59%  "
$ws.Range("D18").Value = "2.619.87"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").Value = "11.43"
$ws.Range("E19").Value = "  +11.14%  "
$ws.Range("E20").Value = "  +3.60%  "
$ws.Range("D21").Value = "348.98"
$ws.Range("E21").Value = "  +3.36%  "
$ws.Range("D22").Value = "6.92"
$ws.Range("E22").Value = "  +7.71%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "0.520"
$ws.Range("E24").Value = "  +9.18%  "
$ws.Range("D25").Value = "63.31"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "7.95"
$ws.Range("E28").Value = "  +7.42%  "
$ws.Range("D29").Value = "0.0₃0799"
$ws.Range("E29").Value = "  +4.23%  "
$ws.Range("E30").Value = "  +11.72%  "
$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "164.31"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("E35").Value = "  +12.80%  "
$ws.Range("D36").Value = "4.26"
$ws.Range("E36").Value = "  +4.21%  "
$ws.Range("E37").Value = "  +6.33%  "
$ws.Range("E38").Value = "  +10.73%  "
$ws.Range("D39").Value = "314.94"
$ws.Range("E39").Value = "  +10.12%  "
$ws.Range("D40").Value = "37.94"
$ws.Range("E40").Value = "  +1.14%  "
$ws.Range("E41").Value = "  +6.38%  "
$ws.Range("D42").Value = "0.843"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "135.25"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "0.0991"
$ws.Range("E44").Value = "  +2.47%  "
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  -0.36%  "
$ws.Range("D46").Value = "19.92"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0554"
$ws.Range("E47").Value = "  +4.75%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "0.608"
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "5.01"
$ws.Range("E49").Value = "  +11.55%  "
$ws.Range("D50").Value = "20.25"
$ws.Range("E50").Value = "  +8.86%  "
$ws.Range("E51").Value = "  +4.57%  "

Write-Host "Updated cryptos list"
